$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$new = $wb.Worksheets.Add($null, $last)

$new.Range("C3").Value = "percent"
$new.Range("D3").Value = "exposure"

$new.Range("C4").Value = 0.33
$new.Range("C4").Style = "Percent"
$new.Range("D4").Value = 2
$new.Range("E4").Formula = "=D4*C4"

$new.Range("H4").Value = 2
$new.Range("I4").Value = 2
$new.Range("J4").Value = 2
$new.Range("K4").Value = 6
$new.Range("L4").Value = 6
$new.Range("M4").Value = 6
$new.Range("N4").Value = 6
$new.Range("O4").Value = 6
$new.Range("P4").Value = 6
$new.Range("Q4").Value = 6

$new.Range("C5").Formula = "=1-C4"
$new.Range("C5").Style = "Percent"
$new.Range("D5").Value = 6
$new.Range("E5").Formula = "=D5*C5"

$new.Range("E6").Formula = "=SUM(E4:E5)"

$new.Range("D10").Formula = "=33.3*2"
$new.Range("G10").Formula = "=3*2"
$new.Range("H10").Formula = "=6*7"
$new.Range("I10").Formula = "=SUM(G10:H10)"

$new.Range("I11").Formula = "=I10/10"

$ws1.Range("H56:I56").Style = "Normal"
$ws1.Range("I69:I71").Clear()

$new.Range("E6").Select()
$new.Activate()
